$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows above the existing first data row (row 2), to hold
# the backward-extended (earlier) real-time GDP data points.
$ws.Rows("2:12").Insert()

# Copy the number formatting / style from the row that now holds the
# former first data row (shifted down to row 13) onto the newly
# inserted rows, so they match the existing date / value formatting.
$ws.Range("A13:B13").Copy()
$ws.Range("A2:B12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the newly inserted rows with the backward-extended data.
$ws.Range("A2").Value = 30681
$ws.Range("B2").Value = 1.466797881812631

$ws.Range("A3").Value = 31047
$ws.Range("B3").Value = 2.900424903011278

$ws.Range("A4").Value = 31412
$ws.Range("B4").Value = 2.603231597845612

$ws.Range("A5").Value = 31777
$ws.Range("B5").Value = 2.279090113735793

$ws.Range("A6").Value = 32142
$ws.Range("B6").Value = 1.278816132757377

$ws.Range("A7").Value = 32508
$ws.Range("B7").Value = 3.44172297297296

$ws.Range("A8").Value = 32873
$ws.Range("B8").Value = 4.033476219636656

$ws.Range("A9").Value = 33238
$ws.Range("B9").Value = 5.482086096613448

$ws.Range("A10").Value = 33603
$ws.Range("B10").Value = 6.118004442050284

$ws.Range("A11").Value = 33969
$ws.Range("B11").Value = 1.839868480884266

$ws.Range("A12").Value = 34334
$ws.Range("B12").Value = -1.202129486518955
